$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G3").Value = "0/3"
$ws.Range("B4").Value = "0/2"
